# Sync attendance_reports, modules_schedules, and assets from main repo - 2025-12-04 18:24:58
#
# The "Recorded By" column (G) stores a comma-separated list of the users/
# systems that recorded/edited a given attendance session. This edit
# normalizes the ordering of those comma-separated names so that the
# literal "System" entry no longer always sorts first.
#
# Concretely (observed from the diff) the following exact-string
# replacements are applied to column G, for every row that contains them:
#
#   "System, backup@backdoor.com, system"  -> "backup@backdoor.com, System, system"
#   "System, dnasr281@gmail.com"           -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com"          -> "backup@backdoor.com, System"
#   "System, admin@admin.com"              -> "admin@admin.com, System"
#   "admin@admin.com, dnasr281@gmail.com"  -> "dnasr281@gmail.com, admin@admin.com"
#
# All other values (e.g. single-author cells, or cells not matching the
# patterns above) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$map = @{
    "System, backup@backdoor.com, system" = "backup@backdoor.com, System, system";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "System, backup@backdoor.com"         = "backup@backdoor.com, System";
    "System, admin@admin.com"             = "admin@admin.com, System";
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com";
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $current = $cell.Value2

    if ($null -ne $current -and $map.ContainsKey($current)) {
        $cell.Value2 = $map[$current]
    }
}
